$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set values in the order that matches the shared-string table append order:
# wanita(137), dp5nt38(138), dp3t24(139), dp3t26(140), jml(141)
$ws.Range("I11").Value = "wanita"
$ws.Range("J1").Value = "dp5nt38"
$ws.Range("K1").Value = "dp3t24"
$ws.Range("P1").Value = "dp3t26"
$ws.Range("L1").Value = "jml"
